$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: objectives text replaced by the docent/teacher identifier
$ws.Cells.Item(10,2).Value = "6376612 - Daisy Rafaela da Silva"
$ws.Cells.Item(10,3).Value = "6376612 - Daisy Rafaela da Silva"

# Row 13 becomes "Programa resumido:" / "Semestral" with a 60pt row height
$ws.Cells.Item(13,1).Value = "Programa resumido:"
$ws.Cells.Item(13,2).Value = "Semestral"
$ws.Cells.Item(13,3).Value = "Semestral"
$ws.Rows.Item(13).RowHeight = 60

# Row 14 becomes just the "Short syllabus:" label (B/C cleared out)
$ws.Cells.Item(14,1).Value = "Short syllabus:"
$ws.Range("B14:C14").Clear()

# Row 15 becomes "Programa:" / "01/01/2012" with a 120pt row height
$ws.Cells.Item(15,1).Value = "Programa:"
$ws.Cells.Item(15,2).Value = "01/01/2012"
$ws.Cells.Item(15,3).Value = "01/01/2012"
$ws.Rows.Item(15).RowHeight = 120

# Row 16 becomes just the "Syllabus:" label (B/C cleared out)
$ws.Cells.Item(16,1).Value = "Syllabus:"
$ws.Range("B16:C16").Clear()

# Row 17 becomes "Avaliação:" with the default (non-custom) row height
$ws.Rows.Item(17).Delete()
$ws.Rows.Item(17).Insert()
$ws.Cells.Item(17,1).Value = "Avaliação:"

# Row 18 becomes "Método:" / docent identifier with a 60pt row height
$ws.Cells.Item(18,1).Value = "Método:"
$ws.Cells.Item(18,2).Value = "6376612 - Daisy Rafaela da Silva"
$ws.Cells.Item(18,3).Value = "6376612 - Daisy Rafaela da Silva"
$ws.Rows.Item(18).RowHeight = 60

# Row 19 label becomes "Critério:" (text stays the same)
$ws.Cells.Item(19,1).Value = "Critério:"

# Row 20 label becomes "Norma de recuperação:" (text stays the same)
$ws.Cells.Item(20,1).Value = "Norma de recuperação:"

# Row 21 becomes "Bibliografia:" with a 120pt row height
$ws.Cells.Item(21,1).Value = "Bibliografia:"
$ws.Rows.Item(21).RowHeight = 120

# Row 22 (old bibliography paragraph) is removed entirely
$ws.Rows.Item(22).Delete()
